$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.0001488876196638067, 0.04240448674262143, 3.900430680208489, 8.660232485948974, 12.60321654051975)
    3  = @(0.6753301551942219, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1.642425054193055)
    4  = @(0.3048080303191223, 1.667794583268128, 3.900430680208489, 0.496779210170732, 6.369812503966472)
    5  = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    6  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    7  = @(1.459612070389937, 1.667794583268128, 3.900430680208489, 0.496779210170732, 7.524616544037286)
    8  = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 3.645393585217082)
    9  = @(0.6753301551942219, 0.04240448674262143, 0.1575252929769615, 8.660232485948974, 9.535492420862779)
    10 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797)
    11 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 14.36450238910742)
    12 = @(1.459612070389937, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 3.074671312995807)
    13 = @(0.01514828764759746, 0.002777888934908601, 3.900430680208489, 0.496779210170732, 4.415136066961727)
    14 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
